$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.42"
$ws.Range("E2").Value = "'1.01%"
$ws.Range("D3").Value = "'36.20"
$ws.Range("E3").Value = "'-2.92%"
$ws.Range("D4").Value = "'5.097"
$ws.Range("E4").Value = "'1.96%"
$ws.Range("D5").Value = "'0.07881"
$ws.Range("E5").Value = "'0.69%"
$ws.Range("D6").Value = "'2.141"
$ws.Range("E6").Value = "'-3.11%"
$ws.Range("D7").Value = "'7.917"
$ws.Range("E7").Value = "'-1.42%"
$ws.Range("D8").Value = "'0.9190"
$ws.Range("E8").Value = "'0.60%"
$ws.Range("D9").Value = "'0.09683"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D10").Value = "'0.1860"
$ws.Range("E10").Value = "'-1.52%"
$ws.Range("D11").Value = "'0.08676"
$ws.Range("E11").Value = "'-0.14%"
$ws.Range("D12").Value = "'0.03553"
$ws.Range("E12").Value = "'0.83%"
$ws.Range("D13").Value = "'0.09930"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("D14").Value = "'0.001442"
$ws.Range("E14").Value = "'-2.65%"
$ws.Range("D15").Value = "'0.005651"
$ws.Range("E15").Value = "'0.22%"
$ws.Range("D16").Value = "'3.442"
$ws.Range("E16").Value = "'-0.49%"
$ws.Range("D17").Value = "'4.105"
$ws.Range("E17").Value = "'1.83%"
$ws.Range("D18").Value = "'2.636"
$ws.Range("E18").Value = "'16.51%"
$ws.Range("D20").Value = "'0.1318"
$ws.Range("E20").Value = "'1.37%"
$ws.Range("D21").Value = "'5.147"
$ws.Range("E21").Value = "'8.12%"
$ws.Range("D22").Value = "'0.2205"
$ws.Range("E22").Value = "'-3.94%"
$ws.Range("D23").Value = "'0.04559"
$ws.Range("E23").Value = "'-1.70%"
$ws.Range("D24").Value = "'0.005056"
$ws.Range("E24").Value = "'5.59%"
$ws.Range("D25").Value = "'0.001235"
$ws.Range("E25").Value = "'0.30%"
$ws.Range("D27").Value = "'0.0004742"
$ws.Range("E27").Value = "'-0.18%"
$ws.Range("D39").Value = "'0.01848"
$ws.Range("E39").Value = "'4.56%"
$ws.Range("D40").Value = "'0.04768"
$ws.Range("E40").Value = "'0.57%"
$ws.Range("D41").Value = "'0.007611"
$ws.Range("E41").Value = "'-5.59%"
$ws.Range("E42").Value = "'0.55%"
$ws.Range("D43").Value = "'0.007721"
$ws.Range("E43").Value = "'0.90%"
$ws.Range("D44").Value = "'0.002183"
$ws.Range("E44").Value = "'-1.26%"
$ws.Range("D45").Value = "'0.01131"
$ws.Range("E45").Value = "'14.61%"
$ws.Range("D46").Value = "'0.00006315"
$ws.Range("E46").Value = "'4.94%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("D48").Value = "'0.0005793"
$ws.Range("E48").Value = "'-0.14%"
$ws.Range("D49").Value = "'47.49"
$ws.Range("E49").Value = "'536.34%"
$ws.Range("D50").Value = "'0.001997"
$ws.Range("E50").Value = "'-25.77%"
$ws.Range("D51").Value = "'0.00002097"
$ws.Range("E51").Value = "'-0.21%"
